$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are plain decimal-looking strings in the source data
# (the sheet mixes thousands-dot-separated and decimal-dot numbers as literal
# text, e.g. "61.018.52" / "529.74"). Writing such a string straight into
# ".Value" lets Excel auto-convert it to a real number, so each cell is first
# switched to the "@" (Text) number format, then written, then has its
# formatting cleared again so it ends up back on the default (unstyled) cell,
# exactly matching every other cell in the sheet.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.744.22'
$ws.Range('D2').ClearFormats()
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.653.18'
$ws.Range('D3').ClearFormats()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '527.52'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.65'
$ws.Range('D6').ClearFormats()
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.582'
$ws.Range('D8').ClearFormats()
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.47'
$ws.Range('D9').ClearFormats()
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.349'
$ws.Range('D11').ClearFormats()
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.115.00'
$ws.Range('D13').ClearFormats()
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '60.745.42'
$ws.Range('D14').ClearFormats()
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.98'
$ws.Range('D15').ClearFormats()
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000142'
$ws.Range('D16').ClearFormats()
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.657.58'
$ws.Range('D17').ClearFormats()
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.76'
$ws.Range('D18').ClearFormats()
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '353.27'
$ws.Range('D19').ClearFormats()
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.64'
$ws.Range('D20').ClearFormats()
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.29'
$ws.Range('D21').ClearFormats()
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.996'
$ws.Range('D22').ClearFormats()
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '61.43'
$ws.Range('D23').ClearFormats()
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.429'
$ws.Range('D24').ClearFormats()
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.167'
$ws.Range('D25').ClearFormats()
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').ClearFormats()
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0₃0851'
$ws.Range('D27').ClearFormats()
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.29'
$ws.Range('D28').ClearFormats()
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.17'
$ws.Range('D30').ClearFormats()
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.42'
$ws.Range('D31').ClearFormats()
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '150.08'
$ws.Range('D33').ClearFormats()
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.11'
$ws.Range('D34').ClearFormats()
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.913'
$ws.Range('D36').ClearFormats()
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.888'
$ws.Range('D37').ClearFormats()
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '36.85'
$ws.Range('D38').ClearFormats()
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '306.30'
$ws.Range('D39').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.47'
$ws.Range('D40').ClearFormats()
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.76'
$ws.Range('D41').ClearFormats()
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.641'
$ws.Range('D42').ClearFormats()
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.101'
$ws.Range('D43').ClearFormats()
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.25'
$ws.Range('D44').ClearFormats()
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0561'
$ws.Range('D45').ClearFormats()
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.997'
$ws.Range('D46').ClearFormats()
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0241'
$ws.Range('D47').ClearFormats()
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.92'
$ws.Range('D48').ClearFormats()
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.22'
$ws.Range('D49').ClearFormats()
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '10.36'
$ws.Range('D50').ClearFormats()
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.989.33'
$ws.Range('D51').ClearFormats()

# Volume(1h) column (E) values already contain non-numeric characters (%, +/-,
# padding spaces) so they are stored as text natively.
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('E3').Value = '  +1.40%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('E5').Value = '  +2.46%  '
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -1.20%  '
$ws.Range('E9').Value = '  -3.69%  '
$ws.Range('E10').Value = '  +4.18%  '
$ws.Range('E11').Value = '  +0.67%  '
$ws.Range('E12').Value = '  -0.59%  '
$ws.Range('E13').Value = '  +1.42%  '
$ws.Range('E14').Value = '  +0.15%  '
$ws.Range('E15').Value = '  +1.12%  '
$ws.Range('E16').Value = '  +0.77%  '
$ws.Range('E17').Value = '  +1.37%  '
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('E19').Value = '  -1.02%  '
$ws.Range('E20').Value = '  -0.40%  '
$ws.Range('E21').Value = '  +1.46%  '
$ws.Range('E22').Value = '  -0.25%  '
$ws.Range('E23').Value = '  +0.98%  '
$ws.Range('E24').Value = '  +0.66%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('E26').Value = '  +0.32%  '
$ws.Range('E27').Value = '  +0.39%  '
$ws.Range('E28').Value = '  -1.26%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('E30').Value = '  +3.07%  '
$ws.Range('E31').Value = '  -0.39%  '
$ws.Range('E32').Value = '  +2.00%  '
$ws.Range('E33').Value = '  -1.17%  '
$ws.Range('E34').Value = '  +1.82%  '
$ws.Range('E35').Value = '  -0.61%  '
$ws.Range('E36').Value = '  +6.72%  '
$ws.Range('E37').Value = '  +0.46%  '
$ws.Range('E38').Value = '  +1.17%  '
$ws.Range('E39').Value = '  +4.31%  '
$ws.Range('E40').Value = '  -1.30%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('E42').Value = '  +2.81%  '
$ws.Range('E43').Value = '  +0.32%  '
$ws.Range('E44').Value = '  +2.28%  '
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('E47').Value = '  +1.94%  '
$ws.Range('E48').Value = '  -0.75%  '
$ws.Range('E49').Value = '  +3.50%  '
$ws.Range('E50').Value = '  +0.59%  '
$ws.Range('E51').Value = '  -0.37%  '
